# Refresh cryptos table (rank links, prices, 1h volume %) to match upstream data pull
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.127.09"
$ws.Range("E2").Value = "  +2.58%  "

# Row 3
$ws.Range("D3").Value = "1.803.46"
$ws.Range("E3").Value = "  +0.64%  "

# Row 4
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'338.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "

# Row 7
$ws.Range("D7").Value = "'0.3923"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.46%  "

# Row 8
$ws.Range("D8").Value = "'0.3482"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.86%  "

# Row 9
$ws.Range("E9").Value = "  -0.95%  "

# Row 10
$ws.Range("D10").Value = "'1.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.22%  "

# Row 11
$ws.Range("D11").Value = "'0.07524"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "

# Row 12
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").Value = "'22.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "

# Row 14
$ws.Range("D14").Value = "'6.497"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "

# Row 15
$ws.Range("D15").Value = "1.809.04"
$ws.Range("E15").Value = "  +1.02%  "

# Row 16
$ws.Range("D16").Value = "'7.131"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.97%  "

# Row 17
$ws.Range("D17").Value = "'0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "

# Row 18
$ws.Range("D18").Value = "'0.06698"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("D19").Value = "'84.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "

# Row 20
$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "

# Row 21
$ws.Range("D21").Value = "'17.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.98%  "

# Row 22
$ws.Range("D22").Value = "'6.543"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "

# Row 23
$ws.Range("D23").Value = "28.113.87"
$ws.Range("E23").Value = "  +2.56%  "

# Row 24
$ws.Range("D24").Value = "'12.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "

# Row 25
$ws.Range("D25").Value = "'2.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.77%  "

# Row 26
$ws.Range("D26").Value = "'1.488"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.73%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'21.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.81%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.510"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.16%  "

# Row 29
$ws.Range("E29").Value = "  +0.20%  "

# Row 30
$ws.Range("D30").Value = "2.014.63"
$ws.Range("E30").Value = "  +1.02%  "

# Row 31
$ws.Range("D31").Value = "'135.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.54%  "

# Row 32
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "'4.026"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.93%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.140"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.78%  "

# Row 34
$ws.Range("D34").Value = "'0.08851"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.91%  "

# Row 35
$ws.Range("D35").Value = "'13.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.33%  "

# Row 36
$ws.Range("D36").Value = "'0.6916"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.37%  "

# Row 37
$ws.Range("D37").Value = "'0.06528"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.39%  "

# Row 38
$ws.Range("B38").Value = "WEMIXTOKEN"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.608"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.94%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02408"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.59%  "

# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'5.417"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.57%  "

# Row 41
$ws.Range("D41").Value = "'0.2207"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

# Row 42
$ws.Range("D42").Value = "'1.253"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.77%  "

# Row 43
$ws.Range("D43").Value = "'8.438"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.75%  "

# Row 44
$ws.Range("E44").Value = "  +1.21%  "

# Row 45
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6401"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.39%  "

# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.869"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.134"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.01%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'130.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.07192"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'79.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "

